$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '59.610.53'
$ws.Range('E2').Value = '  +8.34%  '

# Row 3
$ws.Range('D3').Value = '2.578.70'
$ws.Range('E3').Value = '  +9.98%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '504.86'
$c.ClearFormats()
$ws.Range('E5').Value = '  +6.35%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '155.87'
$c.ClearFormats()
$ws.Range('E6').Value = '  +7.85%  '

# Row 7
$ws.Range('E7').Value = '  +2.83%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').Value = '2.573.47'
$ws.Range('E9').Value = '  +9.76%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '6.13'
$c.ClearFormats()
$ws.Range('E10').Value = '  +13.08%  '

# Row 11
$ws.Range('E11').Value = '  +6.47%  '

# Row 12
$ws.Range('E12').Value = '  +6.14%  '

# Row 13
$ws.Range('E13').Value = '  +1.84%  '

# Row 14
$ws.Range('D14').Value = '3.022.32'
$ws.Range('E14').Value = '  +9.94%  '

# Row 15
$ws.Range('D15').Value = '59.445.55'
$ws.Range('E15').Value = '  +8.06%  '

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '21.81'
$c.ClearFormats()
$ws.Range('E16').Value = '  +9.24%  '

# Row 17
$ws.Range('E17').Value = '  +5.96%  '

# Row 18
$ws.Range('D18').Value = '2.578.50'
$ws.Range('E18').Value = '  +9.84%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.76'
$c.ClearFormats()
$ws.Range('E19').Value = '  +4.35%  '

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '336.59'
$c.ClearFormats()
$ws.Range('E20').Value = '  +6.59%  '

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '10.36'
$c.ClearFormats()
$ws.Range('E21').Value = '  +8.35%  '

# Row 22
$ws.Range('E22').Value = '  +7.96%  '

# Row 23
$ws.Range('E23').Value = '  +0.27%  '

# Row 24
$ws.Range('E24').Value = '  +5.81%  '

# Row 25
$ws.Range('E25').Value = '  +6.14%  '

# Row 26
$ws.Range('E26').Value = '  +7.95%  '

# Row 27
$ws.Range('D27').Value = '2.677.97'
$ws.Range('E27').Value = '  +9.35%  '

# Row 28
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0838'
$ws.Range('E29').Value = '  +12.05%  '

# Row 30
$ws.Range('E30').Value = '  +3.97%  '

# Row 31
$ws.Range('E31').Value = '  +0.01%  '

# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '157.21'
$c.ClearFormats()
$ws.Range('E32').Value = '  +7.50%  '

# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '19.39'
$c.ClearFormats()
$ws.Range('E33').Value = '  +6.63%  '

# Row 34
$ws.Range('E34').Value = '  +6.59%  '

# Row 35
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.52'
$c.ClearFormats()
$ws.Range('E35').Value = '  +9.01%  '

# Row 36
$ws.Range('E36').Value = '  +10.59%  '

# Row 37
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.91'
$c.ClearFormats()
$ws.Range('E37').Value = '  +9.43%  '

# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.851'
$c.ClearFormats()
$ws.Range('E38').Value = '  +5.44%  '

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.75'
$c.ClearFormats()
$ws.Range('E39').Value = '  +10.68%  '

# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '296.74'
$c.ClearFormats()
$ws.Range('E40').Value = '  +18.62%  '

# Row 41
$ws.Range('E41').Value = '  +8.69%  '

# Row 42
$ws.Range('E42').Value = '  +4.70%  '

# Row 43
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.0576'
$c.ClearFormats()
$ws.Range('E43').Value = '  +11.38%  '

# Row 44
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.102'
$c.ClearFormats()
$ws.Range('E44').Value = '  +3.93%  '

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.628'
$c.ClearFormats()
$ws.Range('E45').Value = '  +9.14%  '

# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.784'
$c.ClearFormats()
$ws.Range('E46').Value = '  +26.31%  '

# Row 47
$ws.Range('E47').Value = '  +0.16%  '

# Row 48
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '4.87'
$c.ClearFormats()
$ws.Range('E48').Value = '  +12.89%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '19.14'
$c.ClearFormats()
$ws.Range('E49').Value = '  +15.03%  '

# Row 50
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0236'
$c.ClearFormats()
$ws.Range('E50').Value = '  +7.10%  '

# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '10.25'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.96%  '
